$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 355; existing rows 355-411 shift down to 358-414.
$ws.Rows("355:357").Insert()

# New week block (date 44474) goes into the freshly inserted rows 355-357.
# Column layout: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Tipo,
# G Producto ID, H Producto, I Categoria ID, J Categoria, K Variedad, L Calidad,
# M Volumen, N Precio minimo, O Precio maximo, P Precio promedio ponderado,
# Q Unidad de comercializacion, R Origen, S Precio $/Kg, T Kg/unidad

$rows = @(
    @{ Row = 355; L = "Especial"; M = 440; N = 19000; O = 20000; P = 19500; S = 2786 },
    @{ Row = 356; L = "Primera";  M = 320; N = 14000; O = 15000; P = 14500; S = 2071 },
    @{ Row = 357; L = "Segunda";  M = 260; N = 12000; O = 13000; P = 12500; S = 1786 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44474
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/bandeja 7 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 7
}
